# Update the "want to go" counts (column F) for rows 3-5 on the
# "展览" and "全部类型" worksheets, reflecting newly generated data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 113
    $ws.Range("F4").Value = 142
    $ws.Range("F5").Value = 2995
}
